# Add NSW "Anzac Day (additional)" observance rows for 2026 and 2027.
#
# Anzac Day (25 Apr) falls on a Saturday in 2026 (observed Mon 27 Apr) and on
# a Sunday in 2027 (observed Mon 26 Apr). For those two years we insert an
# extra NSW-only holiday row "Anzac Day (additional)" immediately after the
# existing VIC "Anzac Day" row for that year, shifting all subsequent rows
# down by one (two rows total by the end of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$data = $used.Value2

# Source rows (1-indexed, matching the *original* UsedRange array before any
# insertion) that need a new row inserted directly after them, together with
# the new row's contents. Row 16 = 2026-04-25 VIC Anzac Day; row 43 =
# 2027-04-25 VIC Anzac Day (original numbering, pre-insertion).
$insertAfter = @{
    16 = @("2026-04-27", "NSW", "Anzac Day (additional)")
    43 = @("2027-04-26", "NSW", "Anzac Day (additional)")
}

$newRowCount = $rowCount + $insertAfter.Count
# NOTE: unlike Range.Value2 (1-indexed), a freshly `New-Object`-ed .NET array
# is 0-indexed, so the destination index here runs 0..(newRowCount-1).
$newData = New-Object 'object[,]' $newRowCount, $colCount

$destRow = 0
for ($srcRow = 1; $srcRow -le $rowCount; $srcRow++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $newData[$destRow, ($c - 1)] = $data[$srcRow, $c]
    }
    $destRow++

    if ($insertAfter.ContainsKey($srcRow)) {
        $newRow = $insertAfter[$srcRow]
        for ($c = 1; $c -le $colCount; $c++) {
            $newData[$destRow, ($c - 1)] = $newRow[$c - 1]
        }
        $destRow++
    }
}

# Column A holds dates formatted as plain "yyyy-mm-dd" text (inline strings
# in the original file), not numeric date serials. Force just that data
# column (excluding the header) to Text format before assigning, otherwise
# Excel auto-converts strings that look like dates (e.g. "2026-06-08") into
# date serials. Leaving columns B/C and the header row (r=1) untouched keeps
# the existing bold header style intact.
$dateCol = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($newRowCount, 1))
$dateCol.NumberFormat = "@"

$startCell = $ws.Cells.Item(1, 1)
$endCell = $ws.Cells.Item($newRowCount, $colCount)
$target = $ws.Range($startCell, $endCell)
$target.Value = $newData
